# Applies the edits captured in the commit:
#  - new active selection (G11)
#  - rows 2 & 3 shrink to their auto-fit content height (207 -> 51 / 53)
#    now that the long hyperlink-style wrapped text no longer wraps
#  - the "wrap text" variant of the Hyperlink style is dropped; every
#    hyperlinked cell in rows 2-5 ends up using the plain (non-wrapping)
#    Hyperlink look, except C2/D2/D3 which lose the hyperlink look
#    entirely (back to the default/Normal style)
#  - a new mailto hyperlink is added on D5 (Florian Kleinig's e-mail)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- selection -------------------------------------------------------
$ws.Range("G11").Select()

# --- row heights (auto-fit shrank these once wrapping stopped) -------
$ws.Rows.Item(2).RowHeight = 51
$ws.Rows.Item(3).RowHeight = 53

# --- drop the hyperlink formatting on C2 / D2 / D3 --------------------
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"

# --- remaining hyperlinked cells keep the Hyperlink look, but the
#     wrapText variant of that style is retired, so make sure none of
#     them wrap anymore ------------------------------------------------
$ws.Range("C3").WrapText = $False
$ws.Range("C4").WrapText = $False
$ws.Range("D4").WrapText = $False
$ws.Range("C5").WrapText = $False

# --- add the new hyperlink for D5 (also applies the Hyperlink style) -
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:florian.kleinig@droniq.de")
